{"js": "// The edit replaces the entire single paragraph (which previously held\n// several runs plus a \"_GoBack\" bookmark) with one paragraph containing a\n// single run. That run's text spans many logical \"lines\" joined with\n// plain \"\\n\" characters stored literally inside the <w:t> (not as\n// separate paragraphs / <w:br/> elements), so we build the replacement\n// the same way: clear the body, then insert the whole string at the\n// resulting collapsed range so the embedded newlines survive as literal\n// characters in the text run instead of being promoted to new\n// paragraphs.\nconst body = context.document.body;\nbody.clear();\nawait context.sync();\n\nconst insertionPoint = body.getRange(Word.RangeLocation.start);\nconst newText = \"E7120 \\nPRODUCT TITLE (FULL): Transformers Generations War for Cybertron Optimus Prime\\nPRODUCT TITLE(CONDENSED): Transformers Generations War for Cybertron Optimus Prime\\nPRODUCT RPN: : Transformers Generations Optimus Prime\\nPRODUCT FEATURE BULLETS\\n1 Deluxe Class figure in each wave features an A.I.R. Lock connection point that can attach to Modular Battle Stations\\n1 Deluxe Class figure in each wave converts into a Modular Battle Station mode\\nCompatible with Modular Battle Stations\\nPRODUCT DESCRIPTION:\\nWar for Cybertron Earthrise Deluxe Class figures stand at 5.5 inches and convert into classic G1 Earth modes. They come with their G1-inspired weapon accessories and 1 Deluxe Class figure in each wave can convert into a Modular Battle Station mode. Modular Battle Station modes can connect together in custom configurations to expand and customize the battlefield. Fans can build out epic space battle scenes with these cross-compatible figures! \\nHIDDEN KEYWORDS:\\nLorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum. \\nAssortment #: E7120 (all suffixes)Assortment #: E7124 E7119 E7120 E7121 E7123 (all suffixes)Assortment #: E7124 E7119 E7120 E7121 E7123 (all suffixes)\";\ninsertionPoint.insertText(newText, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# The original single paragraph held several runs plus a \"_GoBack\"\n# bookmark. The target document instead has one paragraph with a single\n# run whose text is the full replacement copy, with the individual\n# logical lines joined by literal \"`n\" (newline) characters stored\n# directly inside the run's text -- not separate paragraphs / line\n# breaks. Clearing the body's content first (which also removes the\n# bookmark) and then assigning the whole string to .Text in one shot\n# reproduces that: the resulting run keeps the embedded newlines as\n# literal characters instead of splitting into new paragraphs.\n$d = $word.ActiveDocument\n$d.Content.Delete()\n$d.Content.Text = \"E7120 `nPRODUCT TITLE (FULL): Transformers Generations War for Cybertron Optimus Prime`nPRODUCT TITLE(CONDENSED): Transformers Generations War for Cybertron Optimus Prime`nPRODUCT RPN: : Transformers Generations Optimus Prime`nPRODUCT FEATURE BULLETS`n1 Deluxe Class figure in each wave features an A.I.R. Lock connection point that can attach to Modular Battle Stations`n1 Deluxe Class figure in each wave converts into a Modular Battle Station mode`nCompatible with Modular Battle Stations`nPRODUCT DESCRIPTION:`nWar for Cybertron Earthrise Deluxe Class figures stand at 5.5 inches and convert into classic G1 Earth modes. They come with their G1-inspired weapon accessories and 1 Deluxe Class figure in each wave can convert into a Modular Battle Station mode. Modular Battle Station modes can connect together in custom configurations to expand and customize the battlefield. Fans can build out epic space battle scenes with these cross-compatible figures! `nHIDDEN KEYWORDS:`nLorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum lorem ipsum. `nAssortment #: E7120 (all suffixes)Assortment #: E7124 E7119 E7120 E7121 E7123 (all suffixes)Assortment #: E7124 E7119 E7120 E7121 E7123 (all suffixes)\"\n"}
